# Logbook 2017-2018 data processing - ending ablation exclusions definition.
# Adds per-row "exclusion window" notes in column F for rows 143-193, and
# marks row 170 (20180423_01_071_072) for deletion with an orange highlight
# and the note "to delete because shit".
#
# Cells are written in top-to-bottom sheet order so new shared-string
# entries land at the same indices the source workbook uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F143").Value = '-1080'
$ws.Range("F145").Value = '7050-'
$ws.Range("F146").Value = '-420'
$ws.Range("F147").Value = '5300-'
$ws.Range("F151").Value = '7560-'
$ws.Range("F152").Value = '-430'
$ws.Range("F156").Value = '3030-'
$ws.Range("F158").Value = '4200-4600 4900-5040 17430-'
$ws.Range("F159").Value = '-100 2545-'
$ws.Range("F160").Value = '-1600'
$ws.Range("F162").Value = '2850-'
$ws.Range("F163").Value = '-740'
$ws.Range("F165").Value = '8900-10600 15915-'
$ws.Range("F166").Value = '-1260'
$ws.Range("F167").Value = '6380-'

# --- Row 170: mark as "to delete" -> orange fill (RGB 255,192,0) across
#     A:D and F:G, plus the note itself in F170. -------------------------
$ws.Range("A170:D170").Interior.Color = 49407
$ws.Range("F170:G170").Interior.Color = 49407
$ws.Range("F170").Value = 'to delete because shit'

$ws.Range("F171").Value = '3090-'
$ws.Range("F172").Value = '-2300'
$ws.Range("F175").Value = '4700-'
$ws.Range("F176").Value = '-1000'
$ws.Range("F177").Value = '5000-'
$ws.Range("F178").Value = '-200'
$ws.Range("F179").Value = '1220-1400 2370-'
$ws.Range("F180").Value = '-390'
$ws.Range("F182").Value = '-1350'
$ws.Range("F184").Value = '6800-'
$ws.Range("F186").Value = '-900 3515-3600'
$ws.Range("F187").Value = '4320-'
$ws.Range("F188").Value = '8100-'
$ws.Range("F189").Value = '-3870 4275-6745 8850-9480 10390-11750'
$ws.Range("F191").Value = '5700-'
$ws.Range("F192").Value = '-1335'
$ws.Range("F193").Value = '2960-'

# --- Update the saved selection to match the end of the edit session ----
$ws.Activate()
$ws.Range("F193").Select()
